$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.41592166666667
$ws.Range("H2").Value = 187.247765
$ws.Range("I2").Value = 0.1654944774607044
$ws.Range("J2").Value = 0.1654944774607044
$ws.Range("M2").Value = 11.25749966666667
$ws.Range("N2").Value = 33.772499
$ws.Range("O2").Value = 0.6929800609896341
$ws.Range("P2").Value = 0.6929800609896341
$ws.Range("Q2").Value = 702.6472173571927
$ws.Range("R2").Value = 6323.824956214735
$ws.Range("S2").Value = 0.1146843730841666
$ws.Range("T2").Value = 0.1146843730841666
$ws.Range("G3").Value = 62.41592166666667
$ws.Range("H3").Value = 187.247765
$ws.Range("I3").Value = 0.1654944774607044
$ws.Range("J3").Value = 0.1654944774607044
$ws.Range("M3").Value = 0.9898276666666668
$ws.Range("O3").Value = 0.06093101107050686
$ws.Range("P3").Value = 0.06093101107050686
$ws.Range("Q3").Value = 61.78100610616612
$ws.Range("R3").Value = 556.0290549554951
$ws.Range("S3").Value = 0.01008374583826593
$ws.Range("T3").Value = 0.01008374583826593
$ws.Range("G4").Value = 62.41592166666667
$ws.Range("H4").Value = 187.247765
$ws.Range("I4").Value = 0.1654944774607044
$ws.Range("J4").Value = 0.1654944774607044
$ws.Range("M4").Value = 3.821582
$ws.Range("N4").Value = 11.464746
$ws.Range("O4").Value = 0.2352458543950409
$ws.Range("P4").Value = 0.2352458543950409
$ws.Range("Q4").Value = 238.5275627547433
$ws.Range("R4").Value = 2146.74806479269
$ws.Range("S4").Value = 0.03893188974790425
$ws.Range("T4").Value = 0.03893188974790426
$ws.Range("G5").Value = 62.41592166666667
$ws.Range("H5").Value = 187.247765
$ws.Range("I5").Value = 0.1654944774607044
$ws.Range("J5").Value = 0.1654944774607044
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1761463333333333
$ws.Range("N5").Value = 0.528439
$ws.Range("O5").Value = 0.01084307354481826
$ws.Range("P5").Value = 0.01084307354481827
$ws.Range("Q5").Value = 10.99433574320389
$ws.Range("R5").Value = 98.949021688835
$ws.Range("S5").Value = 0.001794468790367686
$ws.Range("T5").Value = 0.001794468790367687
$ws.Range("I6").Value = 0.4369365253446571
$ws.Range("J6").Value = 0.436936525344657
$ws.Range("M6").Value = 11.25749966666667
$ws.Range("N6").Value = 33.772499
$ws.Range("O6").Value = 0.6929800609896341
$ws.Range("P6").Value = 0.6929800609896341
$ws.Range("Q6").Value = 1855.120717052579
$ws.Range("R6").Value = 16696.08645347321
$ws.Range("S6").Value = 0.3027882999819392
$ws.Range("T6").Value = 0.3027882999819392
$ws.Range("I7").Value = 0.4369365253446571
$ws.Range("J7").Value = 0.436936525344657
$ws.Range("M7").Value = 0.9898276666666668
$ws.Range("O7").Value = 0.06093101107050686
$ws.Range("P7").Value = 0.06093101107050686
$ws.Range("Q7").Value = 163.1134679206133
$ws.Range("S7").Value = 0.0266229842628841
$ws.Range("T7").Value = 0.0266229842628841
$ws.Range("I8").Value = 0.4369365253446571
$ws.Range("J8").Value = 0.436936525344657
$ws.Range("M8").Value = 3.821582
$ws.Range("N8").Value = 11.464746
$ws.Range("O8").Value = 0.2352458543950409
$ws.Range("P8").Value = 0.2352458543950409
$ws.Range("Q8").Value = 629.7575971605087
$ws.Range("R8").Value = 5667.818374444579
$ws.Range("S8").Value = 0.1027875062211043
$ws.Range("T8").Value = 0.1027875062211043
$ws.Range("I9").Value = 0.4369365253446571
$ws.Range("J9").Value = 0.436936525344657
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1761463333333333
$ws.Range("N9").Value = 0.528439
$ws.Range("O9").Value = 0.01084307354481826
$ws.Range("P9").Value = 0.01084307354481827
$ws.Range("Q9").Value = 29.02711275818078
$ws.Range("R9").Value = 261.244014823627
$ws.Range("S9").Value = 0.004737734878729465
$ws.Range("T9").Value = 0.004737734878729466
$ws.Range("G10").Value = 57.486235
$ws.Range("H10").Value = 172.458705
$ws.Range("I10").Value = 0.1524235190071549
$ws.Range("J10").Value = 0.1524235190071549
$ws.Range("M10").Value = 11.25749966666667
$ws.Range("N10").Value = 33.772499
$ws.Range("O10").Value = 0.6929800609896341
$ws.Range("P10").Value = 0.6929800609896341
$ws.Range("Q10").Value = 647.1512713504217
$ws.Range("R10").Value = 5824.361442153795
$ws.Range("S10").Value = 0.1056264594978329
$ws.Range("T10").Value = 0.1056264594978329
$ws.Range("G11").Value = 57.486235
$ws.Range("H11").Value = 172.458705
$ws.Range("I11").Value = 0.1524235190071549
$ws.Range("J11").Value = 0.1524235190071549
$ws.Range("M11").Value = 0.9898276666666668
$ws.Range("O11").Value = 0.06093101107050686
$ws.Range("P11").Value = 0.06093101107050686
$ws.Range("Q11").Value = 56.90146585550168
$ws.Range("R11").Value = 512.113192699515
$ws.Range("S11").Value = 0.009287319124030568
$ws.Range("T11").Value = 0.009287319124030568
$ws.Range("G12").Value = 57.486235
$ws.Range("H12").Value = 172.458705
$ws.Range("I12").Value = 0.1524235190071549
$ws.Range("J12").Value = 0.1524235190071549
$ws.Range("M12").Value = 3.821582
$ws.Range("N12").Value = 11.464746
$ws.Range("O12").Value = 0.2352458543950409
$ws.Range("P12").Value = 0.2352458543950409
$ws.Range("Q12").Value = 219.68836092377
$ws.Range("R12").Value = 1977.19524831393
$ws.Range("S12").Value = 0.0358570009587369
$ws.Range("T12").Value = 0.03585700095873691
$ws.Range("G13").Value = 57.486235
$ws.Range("H13").Value = 172.458705
$ws.Range("I13").Value = 0.1524235190071549
$ws.Range("J13").Value = 0.1524235190071549
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1761463333333333
$ws.Range("N13").Value = 0.528439
$ws.Range("O13").Value = 0.01084307354481826
$ws.Range("P13").Value = 0.01084307354481827
$ws.Range("Q13").Value = 10.12598951238833
$ws.Range("R13").Value = 91.13390561149501
$ws.Range("S13").Value = 0.001652739426554585
$ws.Range("T13").Value = 0.001652739426554585
$ws.Range("G14").Value = 92.45614233333333
$ws.Range("H14").Value = 277.368427
$ws.Range("I14").Value = 0.2451454781874835
$ws.Range("J14").Value = 0.2451454781874835
$ws.Range("M14").Value = 11.25749966666667
$ws.Range("N14").Value = 33.772499
$ws.Range("O14").Value = 0.6929800609896341
$ws.Range("P14").Value = 0.6929800609896341
$ws.Range("Q14").Value = 1040.824991498786
$ws.Range("R14").Value = 9367.424923489072
$ws.Range("S14").Value = 0.1698809284256953
$ws.Range("T14").Value = 0.1698809284256953
$ws.Range("G15").Value = 92.45614233333333
$ws.Range("H15").Value = 277.368427
$ws.Range("I15").Value = 0.2451454781874835
$ws.Range("J15").Value = 0.2451454781874835
$ws.Range("M15").Value = 0.9898276666666668
$ws.Range("O15").Value = 0.06093101107050686
$ws.Range("P15").Value = 0.06093101107050686
$ws.Range("Q15").Value = 91.51564763480457
$ws.Range("R15").Value = 823.640828713241
$ws.Range("S15").Value = 0.01493696184532626
$ws.Range("T15").Value = 0.01493696184532626
$ws.Range("G16").Value = 92.45614233333333
$ws.Range("H16").Value = 277.368427
$ws.Range("I16").Value = 0.2451454781874835
$ws.Range("J16").Value = 0.2451454781874835
$ws.Range("M16").Value = 3.821582
$ws.Range("N16").Value = 11.464746
$ws.Range("O16").Value = 0.2352458543950409
$ws.Range("P16").Value = 0.2352458543950409
$ws.Range("Q16").Value = 353.3287293305046
$ws.Range("R16").Value = 3179.958563974542
$ws.Range("S16").Value = 0.05766945746729542
$ws.Range("T16").Value = 0.05766945746729542
$ws.Range("G17").Value = 92.45614233333333
$ws.Range("H17").Value = 277.368427
$ws.Range("I17").Value = 0.2451454781874835
$ws.Range("J17").Value = 0.2451454781874835
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1761463333333333
$ws.Range("N17").Value = 0.528439
$ws.Range("O17").Value = 0.01084307354481826
$ws.Range("P17").Value = 0.01084307354481827
$ws.Range("Q17").Value = 16.28581046616144
$ws.Range("R17").Value = 146.572294195453
$ws.Range("S17").Value = 0.002658130449166525
$ws.Range("T17").Value = 0.002658130449166526
